$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1952936666666667
$ws.Range("H2").Value = 0.585881
$ws.Range("I2").Value = 0.001827617096392301
$ws.Range("J2").Value = 0.0018276170963923
$ws.Range("M2").Value = 2.781641666666667
$ws.Range("N2").Value = 8.344925
$ws.Range("O2").Value = 0.1885805187409705
$ws.Range("P2").Value = 0.1885805187409705
$ws.Range("Q2").Value = 0.5432370004361111
$ws.Range("R2").Value = 4.889133003925
$ws.Range("S2").Value = 0.0003446529800975262
$ws.Range("T2").Value = 0.0003446529800975262

$ws.Range("G3").Value = 0.1952936666666667
$ws.Range("H3").Value = 0.585881
$ws.Range("I3").Value = 0.001827617096392301
$ws.Range("J3").Value = 0.0018276170963923
$ws.Range("O3").Value = 0.008234342360037365
$ws.Range("P3").Value = 0.008234342360037365
$ws.Range("Q3").Value = 0.02372036875333333
$ws.Range("R3").Value = 0.21348331878
$ws.Range("S3").Value = [double]"1.504922487475161E-05"
$ws.Range("T3").Value = [double]"1.504922487475161E-05"

$ws.Range("G4").Value = 0.1952936666666667
$ws.Range("H4").Value = 0.585881
$ws.Range("I4").Value = 0.001827617096392301
$ws.Range("J4").Value = 0.0018276170963923
$ws.Range("M4").Value = 5.320086
$ws.Range("N4").Value = 15.960258
$ws.Range("O4").Value = 0.3606735510360756
$ws.Range("P4").Value = 0.3606735510360756
$ws.Range("Q4").Value = 1.038979101922
$ws.Range("R4").Value = 9.350811917298
$ws.Range("S4").Value = 0.0006591731480900528
$ws.Range("T4").Value = 0.0006591731480900526

$ws.Range("G5").Value = 0.1952936666666667
$ws.Range("H5").Value = 0.585881
$ws.Range("I5").Value = 0.001827617096392301
$ws.Range("J5").Value = 0.0018276170963923
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.9776426666666667
$ws.Range("N5").Value = 2.932928
$ws.Range("O5").Value = 0.06627897598479518
$ws.Range("P5").Value = 0.06627897598479518
$ws.Range("Q5").Value = 0.1909274210631111
$ws.Range("R5").Value = 1.718346789568
$ws.Range("S5").Value = 0.0001211325896411864
$ws.Range("T5").Value = 0.0001211325896411864

$ws.Range("G6").Value = 0.1952936666666667
$ws.Range("H6").Value = 0.585881
$ws.Range("I6").Value = 0.001827617096392301
$ws.Range("J6").Value = 0.0018276170963923
$ws.Range("M6").Value = 5.549588666666668
$ws.Range("N6").Value = 16.648766
$ws.Range("O6").Value = 0.3762326118781214
$ws.Range("P6").Value = 0.3762326118781213
$ws.Range("Q6").Value = 1.083799519205111
$ws.Range("R6").Value = 9.754195672846
$ws.Range("S6").Value = 0.0006876091536887836
$ws.Range("T6").Value = 0.0006876091536887833

$ws.Range("G7").Value = 3.363724333333333
$ws.Range("H7").Value = 10.091173
$ws.Range("I7").Value = 0.03147874789838274
$ws.Range("J7").Value = 0.03147874789838274
$ws.Range("M7").Value = 2.781641666666667
$ws.Range("N7").Value = 8.344925
$ws.Range("O7").Value = 0.1885805187409705
$ws.Range("P7").Value = 0.1885805187409705
$ws.Range("Q7").Value = 9.356675760780556
$ws.Range("R7").Value = 84.21008184702499
$ws.Range("S7").Value = 0.005936278607993252
$ws.Range("T7").Value = 0.00593627860799325

$ws.Range("G8").Value = 3.363724333333333
$ws.Range("H8").Value = 10.091173
$ws.Range("I8").Value = 0.03147874789838274
$ws.Range("J8").Value = 0.03147874789838274
$ws.Range("O8").Value = 0.008234342360037365
$ws.Range("P8").Value = 0.008234342360037365
$ws.Range("Q8").Value = 0.4085579575266667
$ws.Range("R8").Value = 3.677021617739999
$ws.Range("S8").Value = 0.0002592067872605902
$ws.Range("T8").Value = 0.0002592067872605902

$ws.Range("G9").Value = 3.363724333333333
$ws.Range("H9").Value = 10.091173
$ws.Range("I9").Value = 0.03147874789838274
$ws.Range("J9").Value = 0.03147874789838274
$ws.Range("M9").Value = 5.320086
$ws.Range("N9").Value = 15.960258
$ws.Range("O9").Value = 0.3606735510360756
$ws.Range("P9").Value = 0.3606735510360756
$ws.Range("Q9").Value = 17.895302733626
$ws.Range("R9").Value = 161.057724602634
$ws.Range("S9").Value = 0.01135355178667911
$ws.Range("T9").Value = 0.01135355178667911

$ws.Range("G10").Value = 3.363724333333333
$ws.Range("H10").Value = 10.091173
$ws.Range("I10").Value = 0.03147874789838274
$ws.Range("J10").Value = 0.03147874789838274
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.9776426666666667
$ws.Range("N10").Value = 2.932928
$ws.Range("O10").Value = 0.06627897598479518
$ws.Range("P10").Value = 0.06627897598479518
$ws.Range("Q10").Value = 3.288520427171556
$ws.Range("R10").Value = 29.596683844544
$ws.Range("S10").Value = 0.002086379175988332
$ws.Range("T10").Value = 0.002086379175988331

$ws.Range("G11").Value = 3.363724333333333
$ws.Range("H11").Value = 10.091173
$ws.Range("I11").Value = 0.03147874789838274
$ws.Range("J11").Value = 0.03147874789838274
$ws.Range("M11").Value = 5.549588666666668
$ws.Range("N11").Value = 16.648766
$ws.Range("O11").Value = 0.3762326118781214
$ws.Range("P11").Value = 0.3762326118781213
$ws.Range("Q11").Value = 18.66728643805756
$ws.Range("R11").Value = 168.005577942518
$ws.Range("S11").Value = 0.01184333154046146
$ws.Range("T11").Value = 0.01184333154046146

$ws.Range("G12").Value = 101.145495
$ws.Range("H12").Value = 303.436485
$ws.Range("I12").Value = 0.9465500804006033
$ws.Range("J12").Value = 0.9465500804006032
$ws.Range("M12").Value = 2.781641666666667
$ws.Range("N12").Value = 8.344925
$ws.Range("O12").Value = 0.1885805187409705
$ws.Range("P12").Value = 0.1885805187409705
$ws.Range("Q12").Value = 281.350523287625
$ws.Range("R12").Value = 2532.154709588625
$ws.Range("S12").Value = 0.1785009051762531
$ws.Range("T12").Value = 0.178500905176253

$ws.Range("G13").Value = 101.145495
$ws.Range("H13").Value = 303.436485
$ws.Range("I13").Value = 0.9465500804006033
$ws.Range("J13").Value = 0.9465500804006032
$ws.Range("O13").Value = 0.008234342360037365
$ws.Range("P13").Value = 0.008234342360037365
$ws.Range("Q13").Value = 12.2851318227
$ws.Range("R13").Value = 110.5661864043
$ws.Range("S13").Value = 0.007794217422939461
$ws.Range("T13").Value = 0.00779421742293946

$ws.Range("G14").Value = 101.145495
$ws.Range("H14").Value = 303.436485
$ws.Range("I14").Value = 0.9465500804006033
$ws.Range("J14").Value = 0.9465500804006032
$ws.Range("M14").Value = 5.320086
$ws.Range("N14").Value = 15.960258
$ws.Range("O14").Value = 0.3606735510360756
$ws.Range("P14").Value = 0.3606735510360756
$ws.Range("Q14").Value = 538.10273191257
$ws.Range("R14").Value = 4842.92458721313
$ws.Range("S14").Value = 0.3413955787315684
$ws.Range("T14").Value = 0.3413955787315684

$ws.Range("G15").Value = 101.145495
$ws.Range("H15").Value = 303.436485
$ws.Range("I15").Value = 0.9465500804006033
$ws.Range("J15").Value = 0.9465500804006032
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.9776426666666667
$ws.Range("N15").Value = 2.932928
$ws.Range("O15").Value = 0.06627897598479518
$ws.Range("P15").Value = 0.06627897598479518
$ws.Range("Q15").Value = 98.88415145312
$ws.Range("R15").Value = 889.95736307808
$ws.Range("S15").Value = 0.06273637004727753
$ws.Range("T15").Value = 0.06273637004727753

$ws.Range("G16").Value = 101.145495
$ws.Range("H16").Value = 303.436485
$ws.Range("I16").Value = 0.9465500804006033
$ws.Range("J16").Value = 0.9465500804006032
$ws.Range("M16").Value = 5.549588666666668
$ws.Range("N16").Value = 16.648766
$ws.Range("O16").Value = 0.3762326118781214
$ws.Range("P16").Value = 0.3762326118781213
$ws.Range("Q16").Value = 561.3158927363901
$ws.Range("R16").Value = 5051.843034627511
$ws.Range("S16").Value = 0.3561230090225648
$ws.Range("T16").Value = 0.3561230090225647

$ws.Range("G17").Value = 2.152479666666667
$ws.Range("H17").Value = 6.457439
$ws.Range("I17").Value = 0.02014355460462176
$ws.Range("J17").Value = 0.02014355460462176
$ws.Range("M17").Value = 2.781641666666667
$ws.Range("N17").Value = 8.344925
$ws.Range("O17").Value = 0.1885805187409705
$ws.Range("P17").Value = 0.1885805187409705
$ws.Range("Q17").Value = 5.987427127452778
$ws.Range("R17").Value = 53.886844147075
$ws.Range("S17").Value = 0.003798681976626635
$ws.Range("T17").Value = 0.003798681976626635

$ws.Range("G18").Value = 2.152479666666667
$ws.Range("H18").Value = 6.457439
$ws.Range("I18").Value = 0.02014355460462176
$ws.Range("J18").Value = 0.02014355460462176
$ws.Range("O18").Value = 0.008234342360037365
$ws.Range("P18").Value = 0.008234342360037365
$ws.Range("Q18").Value = 0.2614401803133333
$ws.Range("R18").Value = 2.35296162282
$ws.Range("S18").Value = 0.0001658689249625627
$ws.Range("T18").Value = 0.0001658689249625626

$ws.Range("G19").Value = 2.152479666666667
$ws.Range("H19").Value = 6.457439
$ws.Range("I19").Value = 0.02014355460462176
$ws.Range("J19").Value = 0.02014355460462176
$ws.Range("M19").Value = 5.320086
$ws.Range("N19").Value = 15.960258
$ws.Range("O19").Value = 0.3606735510360756
$ws.Range("P19").Value = 0.3606735510360756
$ws.Range("Q19").Value = 11.451376939918
$ws.Range("R19").Value = 103.062392459262
$ws.Range("S19").Value = 0.007265247369738022
$ws.Range("T19").Value = 0.007265247369738021

$ws.Range("G20").Value = 2.152479666666667
$ws.Range("H20").Value = 6.457439
$ws.Range("I20").Value = 0.02014355460462176
$ws.Range("J20").Value = 0.02014355460462176
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 0.6666666666666666
$ws.Range("M20").Value = 0.9776426666666667
$ws.Range("N20").Value = 2.932928
$ws.Range("O20").Value = 0.06627897598479518
$ws.Range("P20").Value = 0.06627897598479518
$ws.Range("Q20").Value = 2.104355961265778
$ws.Range("R20").Value = 18.939203651392
$ws.Range("S20").Value = 0.001335094171888136
$ws.Range("T20").Value = 0.001335094171888136

$ws.Range("G21").Value = 2.152479666666667
$ws.Range("H21").Value = 6.457439
$ws.Range("I21").Value = 0.02014355460462176
$ws.Range("J21").Value = 0.02014355460462176
$ws.Range("M21").Value = 5.549588666666668
$ws.Range("N21").Value = 16.648766
$ws.Range("O21").Value = 0.3762326118781214
$ws.Range("P21").Value = 0.3762326118781213
$ws.Range("Q21").Value = 11.94537676336378
$ws.Range("R21").Value = 107.508390870274
$ws.Range("S21").Value = 0.007578662161406404
$ws.Range("T21").Value = 0.007578662161406401
